$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 14:54:06"
$wsZhCn.Range("H2").Value = "2016-03-18 14:54:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 14:54:09"
$wsDeDe.Range("H2").Value = "2016-03-18 14:54:30"
